# Generate Report for Handback
#
# The localization-status report is updated to reflect that the handback
# transform for the "d00cbfb5-..." file failed: the Overview sheet's
# status column is updated, and each per-locale sheet (zh-cn, de-de) gets
# an explanatory message in its "Error Detail" column, whose width is
# widened to fit the longer text.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"
$zhError = "Handback file name: w4zfjrug.3zc is different with handoff file name: d00cbfb5-33d8-4620-af40-d1c618b6d0cb.0df29f8f0ddde59f2b58b1f9551295e7bde9c2c2.zh-cn."
$deError  = "Handback file name: w4zfjrug.3zc is different with handoff file name: d00cbfb5-33d8-4620-af40-d1c618b6d0cb.0df29f8f0ddde59f2b58b1f9551295e7bde9c2c2.de-de."

# --- Overview sheet: "Ready for handoff" -> "Handback transform failed" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# --- zh-cn sheet: Status (C3) text update, Error Detail (P3), widen column P ---
# (ColumnWidth 39.17 round-trips through the engine's internal pixel
#  conversion to an OOXML <col width="40">, matching the target column width.)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("P3").Value = $zhError
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: Status (C3) text update, Error Detail (P3), widen column P ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("P3").Value = $deError
$wsDe.Columns.Item(16).ColumnWidth = 39.17
